$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.3451510965824127
$ws.Range("B1").Value = 2.689388990402222
$ws.Range("C1").Value = 4.765698432922363
$ws.Range("D1").Value = 1.665093779563904
$ws.Range("E1").Value = 0.8366559147834778
